# Refresh crypto price/volume data (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'68.147.94"
$ws.Range('E2').Value = '  +1.55%  '
$ws.Range('D3').Value = "'3.541.73"
$ws.Range('E3').Value = '  +0.74%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'601.95"
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').Value = "'184.81"
$ws.Range('E6').Value = '  +6.62%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = "'0.599"
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range('E9').Value = '  +5.11%  '
$ws.Range('D10').Value = "'7.16"
$ws.Range('E10').Value = '  -1.59%  '
$ws.Range('D11').Value = "'0.447"
$ws.Range('E11').Value = '  +2.32%  '
$ws.Range('D12').Value = "'4.151.42"
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('E13').Value = '  +12.63%  '
$ws.Range('D15').Value = "'68.060.59"
$ws.Range('E15').Value = '  +1.44%  '
$ws.Range('D16').Value = "'0.0000183"
$ws.Range('E16').Value = '  +1.14%  '
$ws.Range('D17').Value = "'3.532.49"
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('D18').Value = "'6.42"
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('D19').Value = "'14.99"
$ws.Range('E19').Value = '  +5.43%  '
$ws.Range('D20').Value = "'401.16"
$ws.Range('E20').Value = '  +1.64%  '
$ws.Range('D21').Value = "'8.18"
$ws.Range('E21').Value = '  +2.26%  '
$ws.Range('D22').Value = "'73.90"
$ws.Range('E22').Value = '  +1.15%  '
$ws.Range('D23').Value = "'0.549"
$ws.Range('E23').Value = '  +1.50%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = "'0.0000126"
$ws.Range('E25').Value = '  +3.22%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = "'5.71"
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').Value = "'10.75"
$ws.Range('E27').Value = '  +5.00%  '
$ws.Range('D28').Value = "'0.180"
$ws.Range('E28').Value = '  -0.67%  '
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = "'6.33"
$ws.Range('E30').Value = '  +1.16%  '
$ws.Range('D31').Value = "'1.47"
$ws.Range('E31').Value = '  +1.51%  '
$ws.Range('D32').Value = "'2.09"
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('D33').Value = "'24.26"
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('D34').Value = "'7.49"
$ws.Range('E34').Value = '  +1.34%  '
$ws.Range('E36').Value = '  +3.32%  '
$ws.Range('D37').Value = "'164.19"
$ws.Range('E37').Value = '  +0.84%  '
$ws.Range('D38').Value = "'1.97"
$ws.Range('E38').Value = '  +2.77%  '
$ws.Range('D39').Value = "'0.885"
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').Value = "'7.19"
$ws.Range('E40').Value = '  +4.07%  '
$ws.Range('D41').Value = "'2.83"
$ws.Range('E41').Value = '  +7.64%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = "'27.47"
$ws.Range('E42').Value = '  +3.92%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').Value = "'4.79"
$ws.Range('E43').Value = '  +2.23%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = "'2.909.94"
$ws.Range('E44').Value = '  +3.75%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = "'27.84"
$ws.Range('E45').Value = '  +0.75%  '
$ws.Range('D46').Value = "'0.0746"
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('D47').Value = "'42.59"
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('D48').Value = "'354.52"
$ws.Range('E48').Value = '  +5.26%  '
$ws.Range('E49').Value = '  +0.76%  '
$ws.Range('E50').Value = '  +0.19%  '
$ws.Range('D51').Value = "'33.90"
$ws.Range('E51').Value = '  +1.11%  '
